# Apply the "st 23. 12. 2020" COVID daily-stats update.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Update AgTests (H) / AgPosit (I) values for existing rows 271-292 ---
# Map: row -> H(new), I(new)  (I = $null means "unchanged")
$updates = @(
    @{ Row = 271; H = 42446; I = 1612 },
    @{ Row = 272; H = 30632; I = 1650 },
    @{ Row = 273; H = 26878; I = 1364 },
    @{ Row = 274; H = 28311; I = 1337 },
    @{ Row = 275; H = 28685; I = 1246 },
    @{ Row = 276; H = 13291; I = 454  },
    @{ Row = 277; H = 3209;  I = $null },
    @{ Row = 279; H = 43610; I = 3099 },
    @{ Row = 280; H = 35941; I = 2392 },
    @{ Row = 281; H = 45611; I = $null },
    @{ Row = 282; H = 46588; I = 2836 },
    @{ Row = 283; H = 17516; I = 1036 },
    @{ Row = 285; H = 40414; I = 3411 },
    @{ Row = 286; H = 54605; I = 4184 },
    @{ Row = 287; H = 56879; I = 3880 },
    @{ Row = 288; H = 53854; I = 3923 },
    @{ Row = 289; H = 62747; I = 3619 },
    @{ Row = 292; H = 76465; I = 6710 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 8).Value = $u.H
    if ($null -ne $u.I) {
        $ws.Cells.Item($u.Row, 9).Value = $u.I
    }
}

# --- 2. Append new row 293 with the full day's data ---
$ws.Cells.Item(293, 1).Value = 44187
$ws.Cells.Item(293, 1).Style = $ws.Cells.Item(292, 1).Style
$ws.Cells.Item(293, 1).NumberFormat = $ws.Cells.Item(292, 1).NumberFormat

$ws.Cells.Item(293, 2).Value = 158905
$ws.Cells.Item(293, 3).Value = 112627
$ws.Cells.Item(293, 4).Value = 44592
$ws.Cells.Item(293, 5).Value = 18144
$ws.Cells.Item(293, 6).Value = 3687
$ws.Cells.Item(293, 7).Value = 1686
$ws.Cells.Item(293, 8).Value = 73961
$ws.Cells.Item(293, 9).Value = 5349
